$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old workbook had 7 trailing species columns (J:O) that are no
# longer measured for this chamber -- drop them (and their header cells)
# for rows 1-5; this also retires their now-orphaned shared strings.
$ws.Range("J1:O5").Clear()

# Row 1 relabeling: this chamber tracks CO2 (new) where "Furan" used to
# sit, and "Furan" (kept) slides into the old "Me-furan" slot; Furfural
# (I1) already reads correctly and every other header is unchanged, so
# leave those alone.
$ws.Range("G1").Value = "CO2"
$ws.Range("H1").Value = "Furan"

# Rows 2-5: only the cells with new calibration numbers get touched.
# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = [double]"6.816250324477081e-07"
$ws.Range("D2").Value = [double]"2.482357144053704e-07"
$ws.Range("I2").Value = [double]"1.393672786043403e-08"

# Row 3
$ws.Range("A3").Value = 0.1
$ws.Range("B3").Value = [double]"5.042037834380973e-06"
$ws.Range("D3").Value = [double]"8.540161573218964e-06"
$ws.Range("I3").Value = [double]"2.500791439450045e-07"

# Row 4
$ws.Range("A4").Value = 0.13
$ws.Range("B4").Value = [double]"3.948567804934743e-06"
$ws.Range("D4").Value = [double]"8.536543819799699e-06"
$ws.Range("I4").Value = [double]"1.688194100535006e-07"

# Row 5
$ws.Range("A5").Value = 0.14
$ws.Range("B5").Value = [double]"4.326178548097709e-06"
$ws.Range("D5").Value = [double]"8.228839267766388e-06"
$ws.Range("I5").Value = [double]"1.629888041917496e-07"

# New calibration rows 6-11.
# Row 6
$ws.Range("A6").Value = 0.15
$ws.Range("B6").Value = [double]"4.416724901953071e-06"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = [double]"1.015107193083668e-05"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = [double]"1.627609185457209e-06"

# Row 7
$ws.Range("A7").Value = 0.16
$ws.Range("B7").Value = [double]"5.984335267943471e-06"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = [double]"1.00470243674512e-05"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = [double]"9.040823630529e-07"

# Row 8
$ws.Range("A8").Value = 0.2
$ws.Range("B8").Value = [double]"4.235344059260807e-06"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = [double]"8.603776128560917e-06"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = [double]"1.48708086255239e-06"

# Row 9
$ws.Range("A9").Value = 0.3
$ws.Range("B9").Value = [double]"5.910972037532766e-06"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = [double]"9.757075148098375e-06"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = [double]"2.242190889275795e-06"

# Row 10
$ws.Range("A10").Value = 0.4
$ws.Range("B10").Value = [double]"6.299399625066266e-06"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = [double]"9.947856412959329e-06"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = [double]"5.601783900947085e-06"

# Row 11
$ws.Range("A11").Value = 0.8
$ws.Range("B11").Value = [double]"5.698042827738062e-06"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = [double]"1.020873319628905e-05"
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = [double]"1.067070616534278e-05"

# Column A carries the bold/bordered/centered header-ish style on every
# row (not just row 1) -- copy it down from A1 onto the freshly added
# rows so A6:A11 match A2:A5.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A6:A11").PasteSpecial(-4122) | Out-Null

